$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5 (shifts existing rows 5-8 down to 6-9)
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5
$ws.Range("A5").Value = "randomLocationForScrollUse"
$ws.Range("B5").Value = "xpath"
$ws.Range("C5").Value = "(//div[contains(@class, 'leaflet-zoom-animated leaflet')]/div)[4]"

# Copy style (border etc.) from row 6 (the row that used to be row 5) into new row 5
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A5:C5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the final row (row 9), which used to be the empty row 8 before insertion
$ws.Range("C9").Value = "//div[@class='leaflet-popup-content']//b[contains(text(), 'Degrees')]"
$ws.Range("A9").Value = "degreesTemperatureOnPopup"
$ws.Range("B9").Value = "xpath"

# A9/B9 pick up a distinct fill-aware style (border + explicit fill) while C9 keeps the plain bordered style
$ws.Range("A9:B9").Interior.Pattern = 1
$ws.Range("A9:B9").Interior.ColorIndex = -4105

# Update selection to the last entry cell, as in the target sheet
$ws.Range("A9").Select()

# Adjust column widths to match new content
$ws.Columns.Item(1).ColumnWidth = 26.453125
$ws.Columns.Item(3).ColumnWidth = 58.36328125
